$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44236
$ws.Range("K2").Value = 'June Pearl'
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 1139

# Row 3
$ws.Range("D3").Value = 44243
$ws.Range("K3").Value = 'Venus'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 21000
$ws.Range("P3").Value = 20500
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("S3").Value = 1139

# Row 4
$ws.Range("D4").Value = 44174
$ws.Range("K4").Value = 'Early John'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20500
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("R4").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S4").Value = 1139

# Row 5
$ws.Range("D5").Value = 44257
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19500
$ws.Range("S5").Value = 1083

# Row 6
$ws.Range("D6").Value = 44223
$ws.Range("K6").Value = 'Ruby Diamond'
$ws.Range("M6").Value = 270
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Región Metropolitana'

# Row 7
$ws.Range("D7").Value = 44223
$ws.Range("K7").Value = 'Super Queen'
$ws.Range("M7").Value = 250
$ws.Range("Q7").Value = '$/caja 18 kilos granel'

# Row 8
$ws.Range("D8").Value = 44209
$ws.Range("K8").Value = 'Super Queen'
$ws.Range("L8").Value = 'Tercera'
$ws.Range("M8").Value = 320
$ws.Range("N8").Value = 17000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 17500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("S8").Value = 972

# Row 9
$ws.Range("D9").Value = 44167
$ws.Range("K9").Value = 'Early John'
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 26000
$ws.Range("P9").Value = 25500
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S9").Value = 1417

# Row 10
$ws.Range("K10").Value = 'Early John'

# Row 11
$ws.Range("D11").Value = 44210
$ws.Range("K11").Value = 'Nectar Crest'
$ws.Range("M11").Value = 250
$ws.Range("Q11").Value = '$/bandeja 18 kilos granel'

# Row 12
$ws.Range("D12").Value = 44210
$ws.Range("K12").Value = 'Red Diamond'
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 19000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 19500
$ws.Range("S12").Value = 1083

# Row 13
$ws.Range("D13").Value = 44161
$ws.Range("K13").Value = 'Artic Glo'
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 280
$ws.Range("N13").Value = 25000
$ws.Range("O13").Value = 26000
$ws.Range("P13").Value = 25500
$ws.Range("S13").Value = 1417

# Row 14
$ws.Range("D14").Value = 44161
$ws.Range("K14").Value = 'Early John'
$ws.Range("N14").Value = 26000
$ws.Range("O14").Value = 27000
$ws.Range("P14").Value = 26500
$ws.Range("Q14").Value = '$/caja 18 kilos granel'
$ws.Range("S14").Value = 1472

# Row 15
$ws.Range("D15").Value = 44238
$ws.Range("K15").Value = 'August Red'
$ws.Range("M15").Value = 320
$ws.Range("R15").Value = 'Región de O''Higgins'

# Row 16
$ws.Range("D16").Value = 44238
$ws.Range("K16").Value = 'Venus'
$ws.Range("M16").Value = 320
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 21000
$ws.Range("P16").Value = 20500
$ws.Range("Q16").Value = '$/bandeja 18 kilos granel'
$ws.Range("S16").Value = 1139

# Row 17
$ws.Range("D17").Value = 44202
$ws.Range("K17").Value = 'Super Queen'
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 20000
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 21000
$ws.Range("S17").Value = 1167

# Row 18
$ws.Range("D18").Value = 44169
$ws.Range("K18").Value = 'Artic Sprite'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 250
$ws.Range("N18").Value = 24000
$ws.Range("O18").Value = 25000
$ws.Range("P18").Value = 24500
$ws.Range("Q18").Value = '$/bandeja 18 kilos granel'
$ws.Range("S18").Value = 1361

# Row 19
$ws.Range("K19").Value = 'Early John'
$ws.Range("M19").Value = 270

# Row 20
$ws.Range("D20").Value = 44216
$ws.Range("K20").Value = 'Nectar Crest'
$ws.Range("M20").Value = 250
$ws.Range("N20").Value = 19000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 19500
$ws.Range("S20").Value = 1083

# Row 21
$ws.Range("D21").Value = 44222
$ws.Range("K21").Value = 'Nectar Crest'
$ws.Range("M21").Value = 270
$ws.Range("N21").Value = 20000
$ws.Range("O21").Value = 21000
$ws.Range("P21").Value = 20500
$ws.Range("S21").Value = 1139

# Row 22
$ws.Range("D22").Value = 44278
$ws.Range("L22").Value = 'Primera'
$ws.Range("N22").Value = 23000
$ws.Range("O22").Value = 24000
$ws.Range("P22").Value = 23500
$ws.Range("S22").Value = 1306

# Row 23
$ws.Range("D23").Value = 44278
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 250
$ws.Range("N23").Value = 23000
$ws.Range("O23").Value = 24000
$ws.Range("P23").Value = 23500
$ws.Range("Q23").Value = '$/bandeja 18 kilos granel'
$ws.Range("S23").Value = 1306

# Row 24
$ws.Range("D24").Value = 44229
$ws.Range("K24").Value = 'Artic Sprite'

# Row 25
$ws.Range("D25").Value = 44201
$ws.Range("K25").Value = 'Super Queen'
$ws.Range("M25").Value = 250
$ws.Range("N25").Value = 22000
$ws.Range("O25").Value = 23000
$ws.Range("P25").Value = 22500
$ws.Range("S25").Value = 1250

# Row 26
$ws.Range("D26").Value = 44215
$ws.Range("K26").Value = 'Venus'
$ws.Range("N26").Value = 19000
$ws.Range("O26").Value = 20000
$ws.Range("P26").Value = 19500
$ws.Range("Q26").Value = '$/bandeja 18 kilos granel'
$ws.Range("S26").Value = 1083

# Row 27
$ws.Range("D27").Value = 44168
$ws.Range("K27").Value = 'Artic Star'
$ws.Range("M27").Value = 270
$ws.Range("N27").Value = 23000
$ws.Range("O27").Value = 24000
$ws.Range("P27").Value = 23500
$ws.Range("Q27").Value = '$/caja 18 kilos granel'
$ws.Range("S27").Value = 1306

# Row 28
$ws.Range("D28").Value = 44168
$ws.Range("K28").Value = 'Early Glo'
$ws.Range("M28").Value = 300
$ws.Range("N28").Value = 23000
$ws.Range("O28").Value = 24000
$ws.Range("P28").Value = 23500
$ws.Range("S28").Value = 1306

# Row 29
$ws.Range("D29").Value = 44217
$ws.Range("K29").Value = 'Ruby Diamond'
$ws.Range("Q29").Value = '$/caja 18 kilos empedrada'

# Row 30
$ws.Range("D30").Value = 44217
$ws.Range("L30").Value = 'Segunda'
$ws.Range("N30").Value = 18000
$ws.Range("O30").Value = 19000
$ws.Range("P30").Value = 18500
$ws.Range("Q30").Value = '$/bandeja 18 kilos granel'
$ws.Range("S30").Value = 1028

# Row 31
$ws.Range("D31").Value = 44244
$ws.Range("K31").Value = 'Nectar Crest'
$ws.Range("N31").Value = 19000
$ws.Range("O31").Value = 20000
$ws.Range("P31").Value = 19500
$ws.Range("S31").Value = 1083

# Row 32
$ws.Range("D32").Value = 44244
$ws.Range("K32").Value = 'Venus'
$ws.Range("M32").Value = 250
$ws.Range("N32").Value = 19000
$ws.Range("O32").Value = 20000
$ws.Range("P32").Value = 19500
$ws.Range("R32").Value = 'Región de O''Higgins'
$ws.Range("S32").Value = 1083

# Row 33
$ws.Range("D33").Value = 44273
$ws.Range("K33").Value = 'Artic Snow'
$ws.Range("N33").Value = 22000
$ws.Range("O33").Value = 23000
$ws.Range("P33").Value = 22500
$ws.Range("Q33").Value = '$/bandeja 18 kilos granel'
$ws.Range("S33").Value = 1250

# Row 34
$ws.Range("D34").Value = 44273
$ws.Range("K34").Value = 'August Red'
$ws.Range("Q34").Value = '$/bandeja 18 kilos granel'
